$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 7787596758.992209
$ws.Range("J4").Value = 7750737057.723518
$ws.Range("K4").Value = 7715246799.71694
$ws.Range("L4").Value = 7795969317.05095
$ws.Range("M4").Value = 7706473785.025129
$ws.Range("N4").Value = 8203978815.343079
$ws.Range("O4").Value = 0.7084717607973422
$ws.Range("P4").Value = 0.7425249169435216
$ws.Range("Q4").Value = 0.7101328903654485
$ws.Range("R4").Value = 0.7735326688815061
$ws.Range("S4").Value = 0.5705980066445183
$ws.Range("T4").Value = 5829.181292952268
$ws.Range("U4").Value = 5718.461575941893
$ws.Range("V4").Value = 5620.234235504641
$ws.Range("W4").Value = 5616.271775782614
$ws.Range("X4").Value = 5489.039181627258
$ws.Range("Y4").Value = 5792.250702123757
$ws.Range("Z4").Value = -110.7197170103755
$ws.Range("AA4").Value = -208.9470574476272
$ws.Range("AB4").Value = -212.9095171696545
$ws.Range("AC4").Value = -340.1421113250108
$ws.Range("AD4").Value = -36.93059082851141
$ws.Range("AE4").Value = -0.01899404246428915
$ws.Range("AF4").Value = -0.03584500926403733
$ws.Range("AG4").Value = -0.03652477191386572
$ws.Range("AH4").Value = -0.05835160963963448
$ws.Range("AI4").Value = -0.006335467876623091
$ws.Range("B5").Value = 3724
$ws.Range("C5").Value = 1319596
$ws.Range("D5").Value = 1340995
$ws.Range("E5").Value = 1363216.5
$ws.Range("F5").Value = 1381377.5
$ws.Range("G5").Value = 1396983.5
$ws.Range("H5").Value = 1405837.5
$ws.Range("I5").Value = 8073698301.538723
$ws.Range("J5").Value = 7973544986.938936
$ws.Range("K5").Value = 7960343364.890577
$ws.Range("L5").Value = 8062966187.272301
$ws.Range("M5").Value = 7938466166.468922
$ws.Range("N5").Value = 8406317901.194114
$ws.Range("O5").Value = 0.7188506981740065
$ws.Range("P5").Value = 0.7526852846401718
$ws.Range("Q5").Value = 0.7231471535982814
$ws.Range("R5").Value = 0.7851772287862513
$ws.Range("S5").Value = 0.5988184747583244
$ws.Range("T5").Value = 6118.31068110143
$ws.Range("U5").Value = 5945.991586052846
$ws.Range("V5").Value = 5839.383080303515
$ws.Range("W5").Value = 5836.902792518556
$ws.Range("X5").Value = 5682.576899776498
$ws.Range("Y5").Value = 5979.580073226183
$ws.Range("Z5").Value = -172.3190950485841
$ws.Range("AA5").Value = -278.9276007979151
$ws.Range("AB5").Value = -281.4078885828749
$ws.Range("AC5").Value = -435.7337813249324
$ws.Range("AD5").Value = -138.7306078752472
$ws.Range("AE5").Value = -0.02816448919157577
$ws.Range("AF5").Value = -0.04558898940184286
$ws.Range("AG5").Value = -0.04599437708387433
$ws.Range("AH5").Value = -0.0712179887613178
$ws.Range("AI5").Value = -0.02267465892239928
$ws.Range("B6").Value = 3388
$ws.Range("C6").Value = 1279930.5
$ws.Range("D6").Value = 1300139
$ws.Range("E6").Value = 1319292.5
$ws.Range("F6").Value = 1335909.5
$ws.Range("G6").Value = 1347445.5
$ws.Range("H6").Value = 1350919.5
$ws.Range("I6").Value = 8195174791.308577
$ws.Range("J6").Value = 8003967668.923697
$ws.Range("K6").Value = 7986190498.832964
$ws.Range("L6").Value = 8053822743.113925
$ws.Range("M6").Value = 7951060294.009595
$ws.Range("N6").Value = 8340224383.655421
$ws.Range("O6").Value = 0.7452774498229043
$ws.Range("P6").Value = 0.7712514757969303
$ws.Range("Q6").Value = 0.7399645808736718
$ws.Range("R6").Value = 0.7966351829988194
$ws.Range("S6").Value = 0.6357733175914995
$ws.Range("T6").Value = 6402.827959259176
$ws.Range("U6").Value = 6156.239962745289
$ws.Range("V6").Value = 6053.388841998999
$ws.Range("W6").Value = 6028.718818987308
$ws.Range("X6").Value = 5900.839992422399
$ws.Range("Y6").Value = 6173.738985672663
$ws.Range("Z6").Value = -246.5879965138874
$ws.Range("AA6").Value = -349.4391172601772
$ws.Range("AB6").Value = -374.1091402718685
$ws.Range("AC6").Value = -501.9879668367776
$ws.Range("AD6").Value = -229.088973586513
$ws.Range("AE6").Value = -0.03851235705268241
$ws.Range("AF6").Value = -0.05457574676121835
$ws.Range("AG6").Value = -0.05842873534199311
$ws.Range("AH6").Value = -0.0784009768856665
$ws.Range("AI6").Value = -0.03577934235375257
$ws.Range("B7").Value = 3498
$ws.Range("C7").Value = 1336284
$ws.Range("D7").Value = 1358919
$ws.Range("E7").Value = 1373575.5
$ws.Range("F7").Value = 1389975
$ws.Range("G7").Value = 1401189.5
$ws.Range("H7").Value = 1401253.5
$ws.Range("I7").Value = 9298898391.511576
$ws.Range("J7").Value = 9052250261.923325
$ws.Range("K7").Value = 9020669333.128187
$ws.Range("L7").Value = 9067888090.133547
$ws.Range("M7").Value = 8931127431.086645
$ws.Range("N7").Value = 9300430238.550495
$ws.Range("O7").Value = 0.7507146941109205
$ws.Range("P7").Value = 0.7607204116638079
$ws.Range("Q7").Value = 0.758147512864494
$ws.Range("R7").Value = 0.8121783876500858
$ws.Range("S7").Value = 0.6558033161806747
$ws.Range("T7").Value = 6958.774026712566
$ws.Range("U7").Value = 6661.361171580737
$ws.Range("V7").Value = 6567.290500688304
$ws.Range("W7").Value = 6523.777830632599
$ws.Range("X7").Value = 6373.961145931115
$ws.Range("Y7").Value = 6637.221772185044
$ws.Range("Z7").Value = -297.412855131829
$ws.Range("AA7").Value = -391.4835260242626
$ws.Range("AB7").Value = -434.9961960799674
$ws.Range("AC7").Value = -584.8128807814519
$ws.Range("AD7").Value = -321.5522545275226
$ws.Range("AE7").Value = -0.04273926039129217
$ws.Range("AF7").Value = -0.05625754256733717
$ws.Range("AG7").Value = -0.06251046440222841
$ws.Range("AH7").Value = -0.08403964240490314
$ws.Range("AI7").Value = -0.04620817593633353
$ws.Range("I8").Value = 9295930975.48064
$ws.Range("J8").Value = 8996298783.328276
$ws.Range("K8").Value = 8935476125.886217
$ws.Range("L8").Value = 8932467611.93041
$ws.Range("M8").Value = 8768962777.859644
$ws.Range("N8").Value = 9064080392.75856
$ws.Range("O8").Value = 0.7500745600954369
$ws.Range("P8").Value = 0.7512675216224277
$ws.Range("Q8").Value = 0.7429167909334924
$ws.Range("R8").Value = 0.7945123769758425
$ws.Range("S8").Value = 0.646585147628989
$ws.Range("T8").Value = 7575.083994322426
$ws.Range("U8").Value = 7229.737684260997
$ws.Range("V8").Value = 7125.709742591943
$ws.Range("W8").Value = 7074.09281760607
$ws.Range("X8").Value = 6927.632656045919
$ws.Range("Y8").Value = 7201.619231870029
$ws.Range("Z8").Value = -345.3463100614299
$ws.Range("AA8").Value = -449.3742517304836
$ws.Range("AB8").Value = -500.9911767163567
$ws.Range("AC8").Value = -647.4513382765072
$ws.Range("AD8").Value = -373.4647624523977
$ws.Range("AE8").Value = -0.04558976643958923
$ws.Range("AF8").Value = -0.05932267577063077
$ws.Range("AG8").Value = -0.06613671572379298
$ws.Range("AH8").Value = -0.08547117612976651
$ws.Range("AI8").Value = -0.0493017321962782
